# Updated AR TCs and Test Plan
$wb = $excel.ActiveWorkbook

# --- CashReceipt: fix mis-spelled "ComapnyID" header to "CompanyID" ---
$wsCashReceipt = $wb.Worksheets.Item("CashReceipt")
$wsCashReceipt.Range("F1").Value = "CompanyID"

# --- CashReceipt_ForeignCurr: insert a new "CompanyID" column before the old column E ---
$wsCashReceiptFC = $wb.Worksheets.Item("CashReceipt_ForeignCurr")
[void]$wsCashReceiptFC.Range("E1").EntireColumn.Insert()
$wsCashReceiptFC.Range("E1").Value = "CompanyID"
$wsCashReceiptFC.Range("E2").Value = "aBb5f0000004JfX"

# --- Update cursor/selection positions on the relevant sheets ---
$wsAddLine = $wb.Worksheets.Item("AddLine")
[void]$wsAddLine.Range("F9").Select()

$wsAddLineFC = $wb.Worksheets.Item("AddLine_ForeignCurr")
[void]$wsAddLineFC.Range("G9").Select()

[void]$wsCashReceipt.Range("F1:F2").Select()

$wsCRATOHome = $wb.Worksheets.Item("CRATO_HomeCurrency")
[void]$wsCRATOHome.Range("C2").Select()

$wsCRATOForeign = $wb.Worksheets.Item("CRATO_ForeignCurrency")
[void]$wsCRATOForeign.Range("F8").Select()

# Select/activate CashReceipt_ForeignCurr last so it becomes the active sheet/tab
[void]$wsCashReceiptFC.Range("H7").Select()
